# Update FuelPrices at 2025-04-14 02:41
# - Append a new row (21) duplicating the last observation (row 20):
#     A21 = 781.86, B21 = 679.38, C21 = 45754 (date serial 2025-04-07)
# - The previously "latest" row (20) keeps its date value (45754) but
#   switches from the "date-only" format to the "date + time" format
#   that every other prior row uses, since it's no longer the newest
#   fetch; the new row 21 takes on the "date-only" formatting instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-format C20: it was the freshest fetch (date-only style) and is now
# superseded, so it gets the standard "date + time" number format used by
# rows 2-19.
$ws.Range("C20").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 21 with the latest fetched values.
$ws.Range("A21").Value = 781.86
$ws.Range("B21").Value = 679.38
$ws.Range("C21").Value = 45754
$ws.Range("C21").NumberFormat = "YYYY-MM-DD"
